$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 220, pushing all existing rows from 220 downward
# (previously rows 220-315) down by one, to rows 221-316.
$ws.Rows.Item(220).Insert()

# Populate the newly inserted row 220 with the new record.
$ws.Range("A220").Value = 1
$ws.Range("B220").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C220").Value = "Arica y Parinacota"
$ws.Range("D220").Value = 44784
$ws.Range("E220").Value = 15
$ws.Range("F220").Value = 100114013
$ws.Range("G220").Value = "Zanahoria"
$ws.Range("H220").Value = "Sin especificar"
$ws.Range("I220").Value = "Primera"
$ws.Range("J220").Value = 80
$ws.Range("K220").Value = 22000
$ws.Range("L220").Value = 23000
$ws.Range("M220").Value = 22500
$ws.Range("N220").Value = '$/saco 25 kilos'
$ws.Range("O220").Value = "Región de Arica y Parinacota"
$ws.Range("P220").Value = 900
$ws.Range("Q220").Value = 25
$ws.Range("R220").Value = "Hortaliza"
